# Fruta / hortaliza, semanal
#
# A new weekly price record was inserted into the "Camote" sheet at row 25,
# pushing the previously existing rows 25-78 down to rows 26-79.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new blank row at position 25 (shifts rows 25..78 down to 26..79,
# and copies the formatting from the row above, same as Excel's native
# Rows.Insert behaviour).
$ws.Rows.Item(25).Insert()

# Populate the newly inserted row 25 with the new record.
$ws.Cells.Item(25, 1).Value = 9
$ws.Cells.Item(25, 2).Value = "Vega Central Mapocho de Santiago"
$ws.Cells.Item(25, 3).Value = "Metropolitana"
$ws.Cells.Item(25, 4).Value = 44662
$ws.Cells.Item(25, 5).Value = 13
$ws.Cells.Item(25, 6).Value = 100114002
$ws.Cells.Item(25, 7).Value = "Camote"
$ws.Cells.Item(25, 8).Value = "Sin especificar"
$ws.Cells.Item(25, 9).Value = "Primera"
$ws.Cells.Item(25, 10).Value = 1600
$ws.Cells.Item(25, 11).Value = 9000
$ws.Cells.Item(25, 12).Value = 10000
$ws.Cells.Item(25, 13).Value = 9500
$ws.Cells.Item(25, 14).Value = '$/malla 18 kilos'
$ws.Cells.Item(25, 15).Value = "Perú"
$ws.Cells.Item(25, 16).Value = 528
$ws.Cells.Item(25, 17).Value = 18
$ws.Cells.Item(25, 18).Value = "Hortaliza"
